# Update crypto price/volume table with refreshed values from GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.322.66"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.692.45"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5404"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.008"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2735"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06457"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07676"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.84%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.542"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.685.40"
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5796"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008409"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.49%  "
$ws.Range("D17").Value = "26.380.37"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.915"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.008"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.268"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "149.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1287"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.865"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06309"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.377"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.327"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.608"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.589"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.676"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.032"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6185"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("E37").Value = "  +1.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01657"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.14%  "
$ws.Range("D39").Value = "1.111.38"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.098"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8849"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.78%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("D44").Value = "1.843.97"
$ws.Range("E44").Value = "  +0.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000113"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.99%  "
$ws.Range("E46").Value = "  +1.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.163"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05282"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4300"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.20%  "
$ws.Range("E51").Value = "  -0.11%  "
